# Invoice-Template.xlsx: add a new "BAG NO." column to the header row of
# Sheet1 (next to the existing "SKU" column), highlight it with a yellow
# fill so it stands out as a newly-added field, and move the active
# selection onto the new column so it is visible when the file is opened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New header cell: AB1 = "BAG NO." (becomes shared-string index 27, right
# after the existing "SKU" entry).
$headerCell = $ws.Range("AB1")
$headerCell.Value = "BAG NO."

# Give the new header a solid yellow fill (RGB 255,255,0 -> FFFFFF00) so it
# is easy to spot; this creates the new fill + cell style the workbook
# needs.
$headerCell.Interior.Color = 65535

# Bring the new column into view and select a cell below it, mirroring the
# reviewer scrolling over to check the freshly added field.
$win = $excel.ActiveWindow
$win.ScrollColumn = 26
$win.ScrollRow = 1
$ws.Range("AB6").Select()

Write-Host "Added BAG NO. column to Sheet1 header."
